$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 250; this shifts the existing rows 250-339
# down to 251-340 and grows the sheet's used range to A1:R340.
$ws.Rows.Item(250).Insert()

# Populate the newly inserted row 250 with the new price-record data.
$ws.Cells.Item(250, 1).Value = 4
$ws.Cells.Item(250, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(250, 3).Value = "Los Lagos"
$ws.Cells.Item(250, 4).Value = 44795
$ws.Cells.Item(250, 5).Value = 10
$ws.Cells.Item(250, 6).Value = 100114014
$ws.Cells.Item(250, 7).Value = "Betarraga"
$ws.Cells.Item(250, 8).Value = "Sin especificar"
$ws.Cells.Item(250, 9).Value = "Primera"
$ws.Cells.Item(250, 10).Value = 750
$ws.Cells.Item(250, 11).Value = 1500
$ws.Cells.Item(250, 12).Value = 1500
$ws.Cells.Item(250, 13).Value = 1500
$ws.Cells.Item(250, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(250, 15).Value = "Región del Maule"
$ws.Cells.Item(250, 16).Value = 300
$ws.Cells.Item(250, 17).Value = 5
$ws.Cells.Item(250, 18).Value = "Hortaliza"
